# Minor typos while recording
# Slide 24: fix "first two solutions" -> "first three solutions" and
# remove the now-obsolete "Can be negative: / no clear meaning" textbox.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(24)

# Fix the caption typo on "TextBox 6"
$s.Shapes.Item("TextBox 6").TextFrame.TextRange.Text = "Wigner function for first three solutions of the harmonic oscillator"

# Remove the "Can be negative: / no clear meaning" textbox entirely
$s.Shapes.Item("TextBox 21").Delete()
